$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Elements sheet: update the canonical terminology URLs ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R200-CanalCommunication/FHIR/TRE-R200-CanalCommunication?vs"
$wsElem.Range("Z5").Value = "https://mos.esante.gouv.fr/NOS/TRE_R256-TypeMessagerie/FHIR/TRE-R256-TypeMessagerie?vs"
$wsElem.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R283-NiveauConfidentialite/FHIR/TRE-R283-NiveauConfidentialite?vs"

# Column Z (26) widened to fit the new, longer URL text.
# (83.333... in COM "characters" units is the closest the pixel-quantized
# ColumnWidth setter can get to the target stored width of 84.1640625.)
$wsElem.Columns.Item(26).ColumnWidth = 83.33333333333333
